$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right after the header row (before existing row 2),
# pushing the current data (old rows 2-5) down to rows 5-8.
# -4121 == xlShiftDown
$ws.Range("A2:D4").Insert(-4121)

# Make sure the newly inserted rows use a plain, non-bold Arial font like
# the rest of the data (Insert() otherwise carries the header's bold font
# down into the new rows).
$ws.Range("A2:D4").Font.Bold = $false
$ws.Range("A2:D4").Font.Name = "Arial"

# Number/date/time formatting for the new rows: Week + Distance stay
# General, Date uses the sheet's date format, Time uses the time format.
$ws.Range("A2:A4").NumberFormat = "General"
$ws.Range("C2:C4").NumberFormat = "General"
$ws.Range("B2:B4").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("D2:D4").NumberFormat = "hh:mm:ss"

# New data for the inserted rows (week 4 entries)
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = 44917
$ws.Range("C2").Value = 2.6
$ws.Range("D2").Value = 0.0162152777777778

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = 44916
$ws.Range("C3").Value = 3.1
$ws.Range("D3").Value = 0.0250347222222222

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = 44915
$ws.Range("C4").Value = 6.2
$ws.Range("D4").Value = 0.0380902777777778

# Renumber the "Week" values of the rows that shifted down (content edit,
# not just a shift) per the target diff.
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 3

$ws.Range("C9").Select()
